{"js": "// Mark the \"\u0110\u00e3 nh\u1eadn C\u00f4ng vi\u1ec7c(X)\" column with an \"X\" for the row whose\n// MSSV (student id) is 0712183 (Ph\u1ea1m Minh Ho\u00e0ng / task \"test\") \u2014 i.e. the\n// student has now accepted (\"nh\u1eadn\") the task, per the commit message.\n\n// The sign-up table is the first table in the document body.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Find the row whose first cell (MSSV) is \"0712183\".\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst targetMssv = \"0712183\";\nlet targetRowIndex = -1;\nfor (let r = 0; r < rows.items.length; r++) {\n  const idCell = table.getCell(r, 0);\n  idCell.load(\"value\");\n  await context.sync();\n  if (idCell.value.trim() === targetMssv) {\n    targetRowIndex = r;\n    break;\n  }\n}\n\nif (targetRowIndex === -1) {\n  throw new Error(\"Could not find row for MSSV \" + targetMssv);\n}\n\n// Last column (index 4) is \"\u0110\u00e3 nh\u1eadn C\u00f4ng vi\u1ec7c(X)\" \u2014 currently empty; type\n// the \"X\" mark there. The cell's paragraph already carries the\n// Times New Roman / sz 26 run formatting via its paragraph mark, so the\n// inserted run picks that formatting up automatically (matching existing\n// rows such as \"0712176\" / \"0712407\" that already have \"X\").\nconst markCell = table.getCell(targetRowIndex, 4);\nmarkCell.body.insertText(\"X\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Mark the \"\u0110\u00e3 nh\u1eadn C\u00f4ng vi\u1ec7c(X)\" column with an \"X\" for the row whose\n# MSSV (student id) is 0712183 (Ph\u1ea1m Minh Ho\u00e0ng / task \"test\") \u2014 i.e. the\n# student has now accepted (\"nh\u1eadn\") the task, per the commit message.\n\n$d = $word.ActiveDocument\n\n# The sign-up table is the first table in the document.\n$t = $d.Tables.Item(1)\n\n$targetMssv = \"0712183\"\n$targetRow = -1\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $idCell = $t.Cell($r, 1)\n    $idText = $idCell.Range.Text -replace \"[\\r\\a]\", \"\"\n    if ($idText.Trim() -eq $targetMssv) {\n        $targetRow = $r\n        break\n    }\n}\n\nif ($targetRow -eq -1) {\n    throw \"Could not find row for MSSV $targetMssv\"\n}\n\n# Last column (5) is \"\u0110\u00e3 nh\u1eadn C\u00f4ng vi\u1ec7c(X)\" \u2014 currently empty; set its text\n# to \"X\". The cell's paragraph mark already carries the Times New Roman /\n# sz 26 run formatting, so the inserted run picks that formatting up\n# automatically (matching existing rows such as \"0712176\" / \"0712407\"\n# that already have \"X\").\n$markCell = $t.Cell($targetRow, 5)\n$markCell.Range.Text = \"X\"\n"}
